$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the D, M, O, P, S values between row 2 and row 4
$cols = @("D", "M", "O", "P", "S")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell4 = $ws.Range($col + "4")
    $v2 = $cell2.Value2
    $v4 = $cell4.Value2
    $cell2.Value2 = $v4
    $cell4.Value2 = $v2
}
